# Update countries & provincias Spain
# - Chile and Corea del Sur swap ranking order (rows 28/29)
# - Pakistan and Rumania swap ranking order (rows 32/33)
# - Refreshed case counts for several countries
# - Updated "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name reorder (ranking swap) ---
$ws.Range("A28").Value = "Chile"
$ws.Range("A29").Value = "Corea del Sur"
$ws.Range("A32").Value = "Pakistan"
$ws.Range("A33").Value = "Rumania"

# --- Updated timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 17:22"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 794330
$ws.Range("C4").Value = 1571
$ws.Range("E4").Value = 679286
$ws.Range("G4").Value = 120
$ws.Range("H4").Value = 42634

# --- Row 9: Reino Unido ---
$ws.Range("B9").Value = 129044
$ws.Range("C9").Value = 4301
$ws.Range("E9").Value = 111363
$ws.Range("G9").Value = 828
$ws.Range("H9").Value = 17337

# --- Row 28: now Chile ---
$ws.Range("B28").Value = 10832
$ws.Range("C28").Value = 325
$ws.Range("D28").Value = 4969
$ws.Range("E28").Value = 5716
$ws.Range("F28").Value = 377
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 147

# --- Row 29: now Corea del Sur ---
$ws.Range("B29").Value = 10683
$ws.Range("C29").Value = 9
$ws.Range("D29").Value = 8213
$ws.Range("E29").Value = 2233
$ws.Range("F29").Value = 55
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 237

# --- Row 32: now Pakistan ---
$ws.Range("B32").Value = 9505
$ws.Range("C32").Value = 613
$ws.Range("D32").Value = 2066
$ws.Range("E32").Value = 7242
$ws.Range("F32").Value = 46
$ws.Range("G32").Value = 21
$ws.Range("H32").Value = 197

# --- Row 33: now Rumania ---
$ws.Range("B33").Value = 9242
$ws.Range("C33").Value = 306
$ws.Range("D33").Value = 2153
$ws.Range("E33").Value = 6607
$ws.Range("F33").Value = 245
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = 482

# --- Row 48: Republica Dominicana ---
$ws.Range("B48").Value = 5044
$ws.Range("C48").Value = 80
$ws.Range("D48").Value = 463
$ws.Range("E48").Value = 4336
$ws.Range("F48").Value = 126
$ws.Range("G48").Value = 10
$ws.Range("H48").Value = 245

# --- Row 64: Kazajistan ---
$ws.Range("D64").Value = 489
$ws.Range("E64").Value = 1469
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 19

# --- Row 85: Bulgaria ---
$ws.Range("B85").Value = 975
$ws.Range("C85").Value = 46
$ws.Range("E85").Value = 760

# --- Row 89: Republica de Chipre ---
$ws.Range("B89").Value = 784
$ws.Range("C89").Value = 12
$ws.Range("D89").Value = 98
$ws.Range("E89").Value = 674

# --- Row 111: Georgia ---
$ws.Range("D111").Value = 97
$ws.Range("E111").Value = 307

# --- Row 114: Montenegro ---
$ws.Range("D114").Value = 101
$ws.Range("E114").Value = 207

# --- Row 157: Bahamas ---
$ws.Range("B157").Value = 64
$ws.Range("C157").Value = 4
$ws.Range("E157").Value = 44
